$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 data: Pedro Henrique -> Carlinhos ---
$ws.Range("C5").Value = "Carlinhos"
$ws.Range("D5").Value = "carlinhos@gmail.com"
$ws.Range("E5").Value = "(48) 99999-9999"

# F5 holds a date-like string ("2005-11-15") that Excel would otherwise
# auto-convert to a date serial. Force it to stay text by temporarily
# stamping a text number format, then restore the original formatting
# (shared with G5) via a format-only copy so the cell's style index is
# unchanged.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2005-11-15"
$ws.Range("G5").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("H5").Value = "MASC"
$ws.Range("I5").Value = "Criciuma"
# G5 (Estudante) and J5 (SC) stay as-is.

# --- Remove the second respondent row entirely ---
$ws.Rows(6).Delete()

# --- Column width adjustments (C, D, E) ---
# Range.ColumnWidth uses character units that render ~0.8333 wider than the
# stored OOXML width attribute for this workbook's default font, so back
# that padding out to land on the exact target widths.
$ws.Columns(3).ColumnWidth = 11 - 0.8333333333333333
$ws.Columns(4).ColumnWidth = 22 - 0.8333333333333333
$ws.Columns(5).ColumnWidth = 18 - 0.8333333333333333

# --- Fix up the active selection to the last remaining data row ---
$ws.Range("B5:J5").Select() | Out-Null
